$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add 5 new survey rows (20-24) below the existing data (which ended at row
# 19), mirroring the formatting already used in the sheet:
#   - Column A: date, formatted like the existing date cells (copy from A19)
#   - Column B: encuestadora name, colour-coded per pollster (copy the fill
#     from an existing cell that already uses that pollster's colour)
#   - Columns C:J: percentages, formatted like the existing percentage cells
# ---------------------------------------------------------------------------

# Seed formatting for each new row by copying an existing, already-styled
# row/cell, then overwrite the values with the new survey data.

# Row 20 - RCN (09-Oct-2023)
$ws.Range("A19").Copy($ws.Range("A20"))
$ws.Range("B10").Copy($ws.Range("B20"))
$ws.Range("C3:J3").Copy($ws.Range("C20:J20"))
$ws.Range("A20").Value = 45208
$ws.Range("B20").Value = "RCN"
$ws.Range("C20").Value = 0.387
$ws.Range("D20").Value = 0.255
$ws.Range("E20").Value = 0.082
$ws.Range("F20").Value = 0.111
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0.16500000000000004

# Row 21 - Mosqueteros (10-Oct-2023)
$ws.Range("A19").Copy($ws.Range("A21"))
$ws.Range("B3").Copy($ws.Range("B21"))
$ws.Range("C3:J3").Copy($ws.Range("C21:J21"))
$ws.Range("A21").Value = 45209
$ws.Range("B21").Value = "Mosqueteros"
$ws.Range("C21").Value = 0.4433
$ws.Range("D21").Value = 0.0792
$ws.Range("E21").Value = 0.0358
$ws.Range("F21").Value = 0.2875
$ws.Range("G21").Value = 0.14
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0.01419999999999999

# Row 22 - RCN (11-Oct-2023)
$ws.Range("A19").Copy($ws.Range("A22"))
$ws.Range("B10").Copy($ws.Range("B22"))
$ws.Range("C3:J3").Copy($ws.Range("C22:J22"))
$ws.Range("A22").Value = 45210
$ws.Range("B22").Value = "RCN"
$ws.Range("C22").Value = 0.371
$ws.Range("D22").Value = 0.272
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0.116
$ws.Range("G22").Value = 0.087
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0.15400000000000003

# Row 23 - RCN (12-Oct-2023)
$ws.Range("A19").Copy($ws.Range("A23"))
$ws.Range("B10").Copy($ws.Range("B23"))
$ws.Range("C3:J3").Copy($ws.Range("C23:J23"))
$ws.Range("A23").Value = 45211
$ws.Range("B23").Value = "RCN"
$ws.Range("C23").Value = 0.356
$ws.Range("D23").Value = 0.288
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0.107
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0.24900000000000011

# Row 24 - Atlas Intel (14-Oct-2023)
$ws.Range("A19").Copy($ws.Range("A24"))
$ws.Range("B14").Copy($ws.Range("B24"))
$ws.Range("C3:J3").Copy($ws.Range("C24:J24"))
$ws.Range("A24").Value = 45213
$ws.Range("B24").Value = "Atlas Intel"
$ws.Range("C24").Value = 0.247
$ws.Range("D24").Value = 0.304
$ws.Range("E24").Value = 0.087
$ws.Range("F24").Value = 0.041
$ws.Range("G24").Value = 0.077
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0.156
$ws.Range("J24").Value = 0.08800000000000008

# Move the selection to reflect where the editor ended up after entering the
# new rows (bottom of the newly-added data).
$ws.Range("L20").Select()
